# v1: treat partial transparency like gray
#
# - remove the now-redundant "make sure auto-detect..." todo (old Id 12)
# - add two new todo items at the top of the Active list:
#     Id 39 "bug: it isn't actually keeping the grayscale"
#     Id 38 "support applying color while zoomed in or out"
# - fix a typo in the "set and check tolerance..." item text
# - bump "Max Id" on the Config sheet from 37 to 39

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Active")

# Helper style reference so new/edited "Created" (date-as-text) cells keep
# the same (default) style as the rest of the table instead of picking up
# an auto-detected date format/style.
$dateStyle = $ws.Cells.Item(2, 5).Style

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $dateStyle
}

# Remove the completed/duplicate item (old Id 12, row 7)
$ws.Rows.Item(7).Delete()

# Insert two fresh rows right after the existing top two entries
$ws.Rows.Item(4).Insert()
$ws.Rows.Item(4).Insert()

# New row 4: Id 39 - "bug: it isn't actually keeping the grayscale"
$ws.Cells.Item(4, 1).Value = 39
$ws.Cells.Item(4, 2).Value = "bug: it isn't actually keeping the grayscale"
$ws.Cells.Item(4, 3).Value = "Todo"
$ws.Cells.Item(4, 4).Value = "Task"
Set-TextValue $ws.Cells.Item(4, 5) "8/12/2018"

# New row 5: Id 38 - "support applying color while zoomed in or out"
$ws.Cells.Item(5, 1).Value = 38
$ws.Cells.Item(5, 2).Value = "support applying color while zoomed in or out"
$ws.Cells.Item(5, 3).Value = "Todo"
$ws.Cells.Item(5, 4).Value = "Task"
Set-TextValue $ws.Cells.Item(5, 5) "8/12/2018"

# Fix a typo ("alter then to white" -> "alter them to white") on the
# tolerance item, now sitting at row 9
$toleranceText = "set and check tolerance for ""black"" and ""white""`n- ""blacks"" will be left untouched`n- ""whites"" will be treated as pure white, which will in effect alter them to white"
$ws.Cells.Item(9, 2).Value = $toleranceText

# Keep the Config sheet's "Max Id" in sync with the new highest Id (39)
$cfg = $wb.Worksheets.Item("Config")
$cfg.Cells.Item(2, 6).Value = 39
